# Auto-update draw results: append the 2025-10-17 Pick 3 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 31

# Columns A, C and E hold strings that look like dates/numbers
# (e.g. "2025-10-17", "251017", an ISO timestamp) in the source data, so force
# text formatting before assigning the values — otherwise Excel would
# auto-coerce them into date serials / numbers instead of keeping them as text.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("C$row").NumberFormat = "@"
$ws.Range("E$row").NumberFormat = "@"

$ws.Range("A$row").Value = "2025-10-17"
$ws.Range("B$row").Value = "Pick 3"
$ws.Range("C$row").Value = "251017"
$ws.Range("D$row").Value = "7-2-8"
$ws.Range("E$row").Value = "2025-10-17T21:37:16.167+04:00"
